# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (fund holdings detail) between the
# existing "总计" (summary) and "2022-Q1" sheets, fills in its data, and
# updates the "总计" sheet with a new summary row for 2022-Q3 (pushing the
# existing 2022-Q1 summary row down one row).

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# --- 1. Insert a new worksheet "2022-Q3" right after "总计" (before "2022-Q1") ---
# NOTE: fetch the "2022-Q1" sheet reference AFTER inserting the new sheet --
# references/lookups made beforehand can rebind to the newly inserted sheet.
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# --- 2. Populate the new "2022-Q3" sheet header row with fund holdings headers ---
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Match header formatting used on the "总计" sheet's header row
$totalSheet.Range("B1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)

# --- 3. Populate the data rows ---
# Numeric-looking text (fund codes with leading zeros, decimal figures) is
# entered with a leading apostrophe to force text storage (otherwise Excel
# auto-converts "013329" -> 13329 and "1.68" -> 1.68 as a number). That
# leaves a "number stored as text" quote-prefix flag on the cell, so we
# immediately paste-format a never-touched default-style cell on top to
# clear the flag again while leaving the text value untouched.
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "'013329"
$q3Sheet.Range("C2").Value = "嘉实全球价值股票（QDII）美元现汇"
$q3Sheet.Range("D2").Value = "'1.68"
$q3Sheet.Range("E2").Value = "'90.63"
$q3Sheet.Range("F2").Value = "'1.24"
$q3Sheet.Range("G2").Value = "'0.0208"
$q3Sheet.Range("H2").Value = 10

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "'013328"
$q3Sheet.Range("C3").Value = "嘉实全球价值股票（QDII）人民币"
$q3Sheet.Range("D3").Value = "'1.68"
$q3Sheet.Range("E3").Value = "'90.63"
$q3Sheet.Range("F3").Value = "'1.24"
$q3Sheet.Range("G3").Value = "'0.0208"
$q3Sheet.Range("H3").Value = 10

# Match the "A" column formatting used on the "总计" sheet
$totalSheet.Range("A2").Copy()
$q3Sheet.Range("A2:A3").PasteSpecial(-4122)

# Clear the quote-prefix flag left behind by the leading-apostrophe entries
# above (copy formatting from an untouched, default-style cell over them).
$q3Sheet.Range("Z99").Copy()
$q3Sheet.Range("B2:B3").PasteSpecial(-4122)
$q3Sheet.Range("D2:G3").PasteSpecial(-4122)

# --- 4. Update the "总计" sheet: push the existing 2022-Q1 summary row down
#        to row 3, then write the new 2022-Q3 summary into row 2 ---
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3:D3").PasteSpecial(-4122)
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3:D3").PasteSpecial(-4163)
$totalSheet.Range("A3").Value = 1

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.04

# Keep "2022-Q1" as the active/selected tab, matching its state before the edit
$q1Sheet.Activate()
